$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")

# Update cell C10 from 18 to 1 (numeric value)
$ws.Range("C10").Value = 1
